{"js": "// Fill in the student's name and DNI on the \"COMPROMISO \u00c9TICO\" placeholder\n// line, underlining the inserted values (thick underline), matching the\n// author's edit:\n//   \"_______________\"  ->  \" Christopher David Pinedo Guti\u00e9rrez \" (underlined)\n//   \"_________\"         ->  \" 72182243 \" (underlined)\n\nconst body = context.document.body;\n\n// --- Replace the blank name line -------------------------------------\nconst nameResults = body.search(\"_______________\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\n\nif (nameResults.items.length > 0) {\n  const nameRange = nameResults.items[0];\n  nameRange.insertText(\" Christopher David Pinedo Guti\u00e9rrez \", \"Replace\");\n  nameRange.font.underline = \"Thick\";\n  await context.sync();\n}\n\n// --- Replace the blank DNI line ---------------------------------------\nconst dniResults = body.search(\"_________\", { matchCase: true });\ndniResults.load(\"items\");\nawait context.sync();\n\nif (dniResults.items.length > 0) {\n  const dniRange = dniResults.items[0];\n  dniRange.insertText(\" 72182243 \", \"Replace\");\n  dniRange.font.underline = \"Thick\";\n  await context.sync();\n}\n", "ps1": "# Fill in the student's name and DNI on the \"COMPROMISO \u00c9TICO\" placeholder\n# line, underlining the inserted values (thick underline), matching the\n# author's edit:\n#   \"_______________\"  ->  \" Christopher David Pinedo Guti\u00e9rrez \" (underlined)\n#   \"_________\"         ->  \" 72182243 \" (underlined)\n\n$d = $word.ActiveDocument\n\n# --- Replace the blank name line -------------------------------------\n$nameRange = $d.Content\n$find = $nameRange.Find\n$find.Text = \"_______________\"\n$find.MatchCase = $true\n$find.MatchWildcards = $false\nif ($find.Execute()) {\n    $nameRange.Text = \" Christopher David Pinedo Guti\u00e9rrez \"\n    $nameRange.Font.Underline = 6\n}\n\n# --- Replace the blank DNI line ---------------------------------------\n$dniRange = $d.Content\n$find2 = $dniRange.Find\n$find2.Text = \"_________\"\n$find2.MatchCase = $true\n$find2.MatchWildcards = $false\nif ($find2.Execute()) {\n    $dniRange.Text = \" 72182243 \"\n    $dniRange.Font.Underline = 6\n}\n"}
